# The product short-name string gets an extra hyphen inserted after "344"
# (e.g. "344MS-EPP-..." -> "344-MS-EPP-...").  That value lives in cell B1
# of both the ProductLoanInput and ProductLoanOutput sheets, so update both.

$wb = $excel.ActiveWorkbook

$wsInput  = $wb.Worksheets.Item("ProductLoanInput")
$wsOutput = $wb.Worksheets.Item("ProductLoanOutput")

$newShortName = "344-MS-EPP-DB-SAR-REC-NON-RNI-CTRFD-DL-MD-TR-1-ONTIME"

$wsInput.Range("B1").Value  = $newShortName
$wsOutput.Range("B1").Value = $newShortName

# Move the selection on the input sheet back up to B1 (it was sitting on B9).
$wsInput.Range("B1").Select()

# The active/selected tab switches from the input sheet to the output sheet,
# whose selection also sits on B1.
$wsOutput.Activate()
$wsOutput.Range("B1").Select()
